$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B with the new simulation output values
$ws.Range("B2").Value = 0.3942291216328126
$ws.Range("B3").Value = 46.36902648081137
$ws.Range("B4").Value = 784.2308943978041
$ws.Range("B5").Value = 73.55231141483992
$ws.Range("B6").Value = 26522.75395529192
$ws.Range("B7").Value = 1686.062464413342
$ws.Range("B8").Value = -1442.708327914368
$ws.Range("B9").Value = 587.3023275902819
$ws.Range("B10").Value = 6707.416755506715
$ws.Range("B11").Value = -521.403569126157
$ws.Range("B12").Value = -8.574458099298266
$ws.Range("B13").Value = -3.921765095485648
$ws.Range("B14").Value = -3.861625337193634
$ws.Range("B15").Value = -0.5503710677441933

# Remove column C entirely (was a second simulation run's results)
$ws.Columns("C").Delete()
